# Updates the per-job Leve profit sheets (currentAveragePrice / LevePrice / LeveProfit
# columns H:N) with refreshed market-board figures, matching the scheduled data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4600
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 4600
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H70").Value = 731418.6
$ws.Range("I70").Value = 3402577
$ws.Range("J70").Value = 2920.9092
$ws.Range("K70").Value = 10207731
$ws.Range("L70").Value = 8762.7276
$ws.Range("M70").Value = -10207461
$ws.Range("N70").Value = -9302.7276
$ws.Range("H73").Value = 731418.6
$ws.Range("I73").Value = 3402577
$ws.Range("J73").Value = 2920.9092
$ws.Range("K73").Value = 10207731
$ws.Range("L73").Value = 8762.7276
$ws.Range("M73").Value = -10206795
$ws.Range("N73").Value = -10634.7276
$ws.Range("H94").Value = 1060.4286
$ws.Range("I94").Value = 1060.4286
$ws.Range("K94").Value = 1060.4286
$ws.Range("M94").Value = -609.4286
$ws.Range("H116").Value = 34380960
$ws.Range("I116").Value = 22824216
$ws.Range("K116").Value = 22824216
$ws.Range("M116").Value = -22820774
$ws.Range("H121").Value = 1785.8226
$ws.Range("J121").Value = 1798.7213
$ws.Range("L121").Value = 5396.1639
$ws.Range("N121").Value = -8890.1639
$ws.Range("H131").Value = 11271.889
$ws.Range("I131").Value = 3670.5715
$ws.Range("K131").Value = 11011.7145
$ws.Range("M131").Value = -5971.7145
$ws.Range("H132").Value = 4422.2964
$ws.Range("I132").Value = 4665.943
$ws.Range("K132").Value = 13997.829
$ws.Range("M132").Value = -11467.829
$ws.Range("H135").Value = 55557212
$ws.Range("I135").Value = 62500616
$ws.Range("K135").Value = 562505544
$ws.Range("M135").Value = -562503009
$ws.Range("H137").Value = 3550.45
$ws.Range("I137").Value = 2187.6667
$ws.Range("J137").Value = 3736.2842
$ws.Range("K137").Value = 6563.000100000001
$ws.Range("L137").Value = 11208.8526
$ws.Range("M137").Value = -4013.000100000001
$ws.Range("N137").Value = -16308.8526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 10509
$ws.Range("I4").Value = 431.66666
$ws.Range("K4").Value = 431.66666
$ws.Range("M4").Value = -315.66666
$ws.Range("H45").Value = 2174.9092
$ws.Range("I45").Value = 1801.375
$ws.Range("K45").Value = 1801.375
$ws.Range("M45").Value = -1424.375
$ws.Range("H54").Value = 44
$ws.Range("I54").Value = 44
$ws.Range("K54").Value = 44
$ws.Range("M54").Value = 725
$ws.Range("H61").Value = 3581.6843
$ws.Range("I61").Value = 2416.3333
$ws.Range("K61").Value = 2416.3333
$ws.Range("M61").Value = -2204.3333
$ws.Range("H74").Value = 4118.864
$ws.Range("I74").Value = 2392.7058
$ws.Range("J74").Value = 9987.799999999999
$ws.Range("K74").Value = 2392.7058
$ws.Range("L74").Value = 9987.799999999999
$ws.Range("M74").Value = -1518.7058
$ws.Range("N74").Value = -11735.8
$ws.Range("H77").Value = 4118.864
$ws.Range("I77").Value = 2392.7058
$ws.Range("J77").Value = 9987.799999999999
$ws.Range("K77").Value = 11963.529
$ws.Range("L77").Value = 49939
$ws.Range("M77").Value = -7595.529
$ws.Range("N77").Value = -58675
$ws.Range("H110").Value = 32264274
$ws.Range("I110").Value = 37038190
$ws.Range("J110").Value = 40320
$ws.Range("K110").Value = 37038190
$ws.Range("L110").Value = 40320
$ws.Range("M110").Value = -37036145
$ws.Range("N110").Value = -44410
$ws.Range("H132").Value = 55559092
$ws.Range("I132").Value = 66670140
$ws.Range("J132").Value = 3861
$ws.Range("K132").Value = 200010420
$ws.Range("L132").Value = 11583
$ws.Range("M132").Value = -200007890
$ws.Range("N132").Value = -16643
$ws.Range("H136").Value = 3581.6843
$ws.Range("I136").Value = 2416.3333
$ws.Range("K136").Value = 7248.999899999999
$ws.Range("M136").Value = -4698.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 325
$ws.Range("I22").Value = 325
$ws.Range("K22").Value = 325
$ws.Range("M22").Value = -152
$ws.Range("H86").Value = 27780160
$ws.Range("I86").Value = 41668540
$ws.Range("K86").Value = 41668540
$ws.Range("M86").Value = -41667417
$ws.Range("H89").Value = 27780160
$ws.Range("I89").Value = 41668540
$ws.Range("K89").Value = 208342700
$ws.Range("M89").Value = -208337084
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H105").Value = 1219.8636
$ws.Range("I105").Value = 979.05884
$ws.Range("K105").Value = 979.05884
$ws.Range("M105").Value = 767.94116
$ws.Range("H107").Value = 21766198
$ws.Range("I107").Value = 13783.333
$ws.Range("K107").Value = 13783.333
$ws.Range("M107").Value = -11863.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3504.0967
$ws.Range("I31").Value = 1050.6818
$ws.Range("J31").Value = 4264.31
$ws.Range("K31").Value = 1050.6818
$ws.Range("L31").Value = 4264.31
$ws.Range("M31").Value = -755.6818000000001
$ws.Range("N31").Value = -4854.31
$ws.Range("H34").Value = 3504.0967
$ws.Range("I34").Value = 1050.6818
$ws.Range("J34").Value = 4264.31
$ws.Range("K34").Value = 1050.6818
$ws.Range("L34").Value = 4264.31
$ws.Range("M34").Value = -848.6818000000001
$ws.Range("N34").Value = -4668.31
$ws.Range("H96").Value = 35000
$ws.Range("J96").Value = 35000
$ws.Range("L96").Value = 35000
$ws.Range("N96").Value = -40492
$ws.Range("H105").Value = 1996.6666
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 1990
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 1990
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -5484
$ws.Range("H132").Value = 3337166.2
$ws.Range("I132").Value = 5003999.5
$ws.Range("J132").Value = 2503749.8
$ws.Range("K132").Value = 15011998.5
$ws.Range("L132").Value = 7511249.399999999
$ws.Range("M132").Value = -15009468.5
$ws.Range("N132").Value = -7516309.399999999
$ws.Range("H134").Value = 2155.318
$ws.Range("I134").Value = 1369.875
$ws.Range("J134").Value = 4249.8335
$ws.Range("K134").Value = 4109.625
$ws.Range("L134").Value = 12749.5005
$ws.Range("M134").Value = -1574.625
$ws.Range("N134").Value = -17819.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 13157945
$ws.Range("I2").Value = 20.545454
$ws.Range("K2").Value = 123.272724
$ws.Range("M2").Value = -10.272724
$ws.Range("H4").Value = 10735011
$ws.Range("I4").Value = 17338746
$ws.Range("K4").Value = 52016238
$ws.Range("M4").Value = -52016126
$ws.Range("H23").Value = 4334.3335
$ws.Range("J23").Value = 5001
$ws.Range("L23").Value = 15003
$ws.Range("N23").Value = -15473
$ws.Range("H68").Value = 2099.8462
$ws.Range("I68").Value = 1362.25
$ws.Range("J68").Value = 3280
$ws.Range("K68").Value = 4086.75
$ws.Range("L68").Value = 9840
$ws.Range("M68").Value = -3275.75
$ws.Range("N68").Value = -11462
$ws.Range("H71").Value = 2099.8462
$ws.Range("I71").Value = 1362.25
$ws.Range("J71").Value = 3280
$ws.Range("K71").Value = 12260.25
$ws.Range("L71").Value = 29520
$ws.Range("M71").Value = -8204.25
$ws.Range("N71").Value = -37632
$ws.Range("H107").Value = 1608.5807
$ws.Range("I107").Value = 1428.0588
$ws.Range("K107").Value = 4284.1764
$ws.Range("M107").Value = -2364.1764

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2561.8096
$ws.Range("I102").Value = 1072.6364
$ws.Range("K102").Value = 1072.6364
$ws.Range("M102").Value = 549.3635999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15627292
$ws.Range("I7").Value = 19232966
$ws.Range("K7").Value = 19232966
$ws.Range("M7").Value = -19232854
$ws.Range("H40").Value = 2723
$ws.Range("I40").Value = 2449.4827
$ws.Range("J40").Value = 3856.1428
$ws.Range("K40").Value = 2449.4827
$ws.Range("L40").Value = 3856.1428
$ws.Range("M40").Value = -2313.4827
$ws.Range("N40").Value = -4128.1428
$ws.Range("H93").Value = 2561.3914
$ws.Range("I93").Value = 2762.2666
$ws.Range("J93").Value = 2184.75
$ws.Range("K93").Value = 2762.2666
$ws.Range("L93").Value = 2184.75
$ws.Range("M93").Value = -1514.2666
$ws.Range("N93").Value = -4680.75
$ws.Range("H122").Value = 4576.1577
$ws.Range("I122").Value = 3883.111
$ws.Range("K122").Value = 11649.333
$ws.Range("M122").Value = -9199.332999999999
$ws.Range("H126").Value = 15627292
$ws.Range("I126").Value = 19232966
$ws.Range("K126").Value = 57698898
$ws.Range("M126").Value = -57696428
$ws.Range("H136").Value = 7095.387
$ws.Range("I136").Value = 6854
$ws.Range("K136").Value = 20562
$ws.Range("M136").Value = -18012

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 5596.7
$ws.Range("I100").Value = 7892
$ws.Range("J100").Value = 241
$ws.Range("K100").Value = 15784
$ws.Range("L100").Value = 482
$ws.Range("M100").Value = -15243
$ws.Range("N100").Value = -1564
$ws.Range("H122").Value = 2250.3635
$ws.Range("J122").Value = 2248.3333
$ws.Range("L122").Value = 6744.999899999999
$ws.Range("N122").Value = -11644.9999
$ws.Range("H132").Value = 838982.8
$ws.Range("I132").Value = 2507106.5
$ws.Range("J132").Value = 4920.875
$ws.Range("K132").Value = 7521319.5
$ws.Range("L132").Value = 14762.625
$ws.Range("M132").Value = -7518789.5
$ws.Range("N132").Value = -19822.625
$ws.Range("H136").Value = 2673.6445
$ws.Range("I136").Value = 1968.5
$ws.Range("K136").Value = 5905.5
$ws.Range("M136").Value = -3355.5
